$wb = $excel.ActiveWorkbook

$wsMeta  = $wb.Worksheets.Item("Meta")
$wsInstr = $wb.Worksheets.Item("Instructions")

# --- Text content updates -------------------------------------------------

# Instructions!B6 - drop "PON" from the list of sample-fraction-distinguished
# parameters (shared string edited in place, keeps its original index).
$wsInstr.Range("B6").Value = "Name of the measured parameter.`n- Note that this can be either the WQX or Simple parameter name.  However, if a parameter is distinguished by Sample Fraction only (i.e. TDP, TDN), then the Simple parameter name must be used here and in all other files (Results, DQO, etc.)"

# Instructions!D9 - more detailed "Available Values" note for Result Sample
# Fraction (this needs to be created before the "MassBays" string below so
# new shared strings land in the same order as the target workbook).
$wsInstr.Range("D9").Value = "standard list in WQX`nRecommended:  Filtered, lab; Filtered, field; Unfiltered; Non-Filterable (Particle)"
$wsInstr.Range("D9").WrapText = $true
$wsInstr.Rows.Item(9).RowHeight = 45

# Meta!B4:B6 - example "Sampling Method Context" switched from MassWateR to
# MassBays.
$wsMeta.Range("B4").Value = "MassBays"
$wsMeta.Range("B5").Value = "MassBays"
$wsMeta.Range("B6").Value = "MassBays"

# --- Alignment tweaks on the Meta sheet ------------------------------------
# Column A (Parameter header + values) becomes explicitly left aligned.
$wsMeta.Range("A1").HorizontalAlignment = -4131   # xlLeft
$wsMeta.Range("A2:A6").HorizontalAlignment = -4131   # xlLeft

# --- Selection / active sheet state ----------------------------------------
# Instructions tab loses the "active" marker and its remembered selection
# moves to D10.
$wsInstr.Activate()
$null = $wsInstr.Range("D10").Select()

# Meta becomes the active/selected tab, with B10 remembered as the selection.
$wsMeta.Activate()
$null = $wsMeta.Range("B10").Select()
